$d = $word.ActiveDocument

# --- 1. Merge the split "Milestone " + "<N>" runs into a single run, for
#        Milestones 3 through 13 (Milestones 1 and 2 are already single
#        runs; Milestone 12 is split across three runs: "Milestone ", "1",
#        "2"). A Find/Replace of the full visible label collapses the
#        paragraph's runs into one, matching the target structure.
foreach ($n in 3..13) {
    $label = "Milestone $n"
    $d.Content.Find.Execute($label, $true, $false, $false, $false, $false, `
        $true, 1, $false, $label, 2) | Out-Null
}

# --- 2. Find the "Maintenance/Prep" paragraph that belongs to Milestone 12
#        (the one immediately following the "Milestone 12" paragraph --
#        Milestones 10, 11 and 13 each have their own separate
#        "Maintenance/Prep" entry that must stay untouched).
$paraCount = $d.Paragraphs.Count
$maintIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "Milestone 12") {
        $maintIdx = $i + 1
        break
    }
}

$maintRange = $d.Paragraphs.Item($maintIdx).Range
# Exclude the paragraph mark so the search/insert stays inside this cell.
$maintRange.MoveEnd(1, -1) | Out-Null
$maintRange.Collapse(0)
$maintRange.InsertAfter("/Tutorial")

# --- 3. Move the "_GoBack" bookmark from the final (empty) paragraph of the
#        document body to the end of the paragraph just edited above --
#        standard Word behaviour is for this bookmark to track the most
#        recent edit location.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$maintRange2 = $d.Paragraphs.Item($maintIdx).Range
$maintRange2.MoveEnd(1, -1) | Out-Null
$maintRange2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $maintRange2) | Out-Null

# --- 4. Merge the header's five text runs ("AR Pets", spacer, "Project
#        Plan", spacer, "2/4/2019") into a single run. Replacing a short
#        anchor substring with itself is enough to make the whole
#        paragraph's runs collapse into one.
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrPara = $hdr.Range.Paragraphs.Item(1)
$hdrRange = $hdrPara.Range
$hdrRange.Find.Execute("AR Pets", $true, $false, $false, $false, $false, `
    $true, 1, $false, "AR Pets", 2) | Out-Null
